$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1030.9
$ws.Range("I38").Value = 1030.9
$ws.Range("K38").Value = 3092.7
$ws.Range("M38").Value = -2720.7
$ws.Range("H40").Value = 3688
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3688
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 3688
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -4038

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4556.048
$ws.Range("I32").Value = 4556.048
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4556.048
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4269.048
$ws.Range("N32").ClearContents()
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H132").Value = 4596.8887
$ws.Range("I132").Value = 4702.8184
$ws.Range("J132").Value = 4430.4287
$ws.Range("K132").Value = 14108.4552
$ws.Range("L132").Value = 13291.2861
$ws.Range("M132").Value = -11578.4552
$ws.Range("N132").Value = -18351.2861
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7192.2
$ws.Range("I134").Value = 7192.2
$ws.Range("K134").Value = 21576.6
$ws.Range("M134").Value = -19041.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
$ws.Range("H74").Value = 39997.727
$ws.Range("J74").Value = 39997.727
$ws.Range("L74").Value = 39997.727
$ws.Range("N74").Value = -41745.727
$ws.Range("H77").Value = 39997.727
$ws.Range("J77").Value = 39997.727
$ws.Range("L77").Value = 119993.181
$ws.Range("N77").Value = -128729.181
$ws.Range("H81").Value = 34666.668
$ws.Range("J81").Value = 34666.668
$ws.Range("L81").Value = 34666.668
$ws.Range("N81").Value = -36662.668
$ws.Range("H82").Value = 22500
$ws.Range("J82").Value = 30000
$ws.Range("L82").Value = 30000
$ws.Range("N82").Value = -30722
$ws.Range("H84").Value = 34666.668
$ws.Range("J84").Value = 34666.668
$ws.Range("L84").Value = 104000.004
$ws.Range("N84").Value = -113984.004
$ws.Range("H85").Value = 22500
$ws.Range("J85").Value = 30000
$ws.Range("L85").Value = 30000
$ws.Range("N85").Value = -32496
$ws.Range("H88").Value = 10198.5
$ws.Range("J88").Value = 10198.5
$ws.Range("L88").Value = 10198.5
$ws.Range("N88").Value = -11010.5
$ws.Range("H91").Value = 10198.5
$ws.Range("J91").Value = 10198.5
$ws.Range("L91").Value = 10198.5
$ws.Range("N91").Value = -13006.5
$ws.Range("H99").Value = 7950
$ws.Range("I99").Value = 8100
$ws.Range("J99").Value = 7500
$ws.Range("K99").Value = 8100
$ws.Range("L99").Value = 7500
$ws.Range("M99").Value = -6602
$ws.Range("N99").Value = -10496
$ws.Range("H126").Value = 7950
$ws.Range("I126").Value = 8100
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 24300
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -21830
$ws.Range("N126").Value = -27440
$ws.Range("H132").Value = 3250
$ws.Range("J132").Value = 750
$ws.Range("L132").Value = 2250
$ws.Range("N132").Value = -7310
$ws.Range("H134").Value = 1656.8334
$ws.Range("I134").Value = 1587.6818
$ws.Range("J134").Value = 2417.5
$ws.Range("K134").Value = 4763.0454
$ws.Range("L134").Value = 7252.5
$ws.Range("M134").Value = -2228.0454
$ws.Range("N134").Value = -12322.5
$ws.Range("H141").Value = 100000
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 227
$ws.Range("I2").Value = 195.3077
$ws.Range("J2").Value = 330
$ws.Range("K2").Value = 1171.8462
$ws.Range("L2").Value = 1980
$ws.Range("M2").Value = -1058.8462
$ws.Range("N2").Value = -2206
$ws.Range("H12").Value = 187.07143
$ws.Range("J12").Value = 214
$ws.Range("L12").Value = 642
$ws.Range("N12").Value = -988
$ws.Range("H38").Value = 138.22223
$ws.Range("I38").Value = 188
$ws.Range("J38").Value = 38.666668
$ws.Range("K38").Value = 564
$ws.Range("L38").Value = 116.000004
$ws.Range("M38").Value = -217
$ws.Range("N38").Value = -810.000004
$ws.Range("H128").Value = 344666.66
$ws.Range("I128").Value = 344666.66
$ws.Range("K128").Value = 1033999.98
$ws.Range("M128").Value = -1029019.98

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 883.64703
$ws.Range("I97").Value = 876.4375
$ws.Range("K97").Value = 876.4375
$ws.Range("M97").Value = -380.4375
$ws.Range("H122").Value = 21863.354
$ws.Range("J122").Value = 35249.75
$ws.Range("L122").Value = 105749.25
$ws.Range("N122").Value = -110649.25
$ws.Range("H126").Value = 1200
$ws.Range("I126").Value = 1200
$ws.Range("K126").Value = 3600
$ws.Range("M126").Value = -1130

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 1266.6666
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 1800
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 1800
$ws.Range("M43").Value = -807
$ws.Range("N43").Value = -2186
$ws.Range("H46").Value = 4340.6665
$ws.Range("I46").Value = 4544.4443
$ws.Range("K46").Value = 4544.4443
$ws.Range("M46").Value = -4356.4443
$ws.Range("H69").Value = 55850
$ws.Range("J69").Value = 55850
$ws.Range("L69").Value = 55850
$ws.Range("N69").Value = -57472
$ws.Range("H72").Value = 55850
$ws.Range("J72").Value = 55850
$ws.Range("L72").Value = 167550
$ws.Range("N72").Value = -175662

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2058.4
$ws.Range("I126").Value = 2058.4
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6175.200000000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3705.200000000001
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 1995
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1995
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 5985
$ws.Range("N132").Value = -11045
$ws.Range("H135").Value = 546870.4
$ws.Range("J135").Value = 56213.125
$ws.Range("L135").Value = 56213.125
$ws.Range("N135").Value = -66353.125
